# Apply the latest cryptos snapshot: updated prices/volumes and the
# Dai/BitcoinCash row swap (rows 19-20), per the Sep 14 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks numeric need NumberFormat "@" so Excel keeps
# them as text (matching the source data) instead of coercing to a number;
# the Style reset afterwards keeps the cell on the default (unstyled) xf so
# no stray style gets attached to the cell.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '26.712.95'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '1.636.05'
$ws.Range("E4").Value = '  +0.08%  '
Set-TextValue "D5" '213.31'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("E8").Value = '  +0.52%  '
$ws.Range("E9").Value = '  +0.69%  '
Set-TextValue "D10" '19.10'
$ws.Range("E10").Value = '  +3.02%  '
Set-TextValue "D11" '0.0836'
$ws.Range("E11").Value = '  +2.64%  '
$ws.Range("D12").Value = '1.862.58'
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("D13").Value = '1.643.26'
$ws.Range("E13").Value = '  +1.70%  '
Set-TextValue "D15" '0.524'
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("D16").Value = '26.673.15'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D19" '209.48'
$ws.Range("E19").Value = '  +3.10%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D20" '1.00'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("E21").Value = '  +0.50%  '
Set-TextValue "D22" '9.37'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("E24").Value = '  +0.57%  '
Set-TextValue "D25" '146.33'
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("E28").Value = '  +0.81%  '
Set-TextValue "D29" '6.68'
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("E30").Value = '  +5.46%  '
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("E33").Value = '  -0.04%  '
Set-TextValue "D34" '1.50'
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("D36").Value = '1.164.53'
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("E37").Value = '  +0.91%  '
Set-TextValue "D38" '0.808'
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("E40").Value = '  -0.01%  '
Set-TextValue "D41" '0.501'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("E43").Value = '  +2.02%  '
$ws.Range("D44").Value = '1.772.98'
$ws.Range("E44").Value = '  +1.16%  '
Set-TextValue "D45" '92.58'
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("D47").Value = '0.0₆0104'
$ws.Range("E47").Value = '  +6.82%  '
Set-TextValue "D48" '54.60'
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("E50").Value = '  +4.46%  '
Set-TextValue "D51" '0.410'
$ws.Range("E51").Value = '  +0.83%  '
